$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.147.86'
$ws.Range("E2").Value = '  -6.10%  '
$ws.Range("D3").Value = '2.552.58'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '299.18'
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("D6").Value = '94.38'
$ws.Range("E6").Value = '  -4.53%  '
$ws.Range("D7").Value = '0.575'
$ws.Range("E7").Value = '  -3.07%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.550'
$ws.Range("E9").Value = '  -4.86%  '
$ws.Range("D10").Value = '36.17'
$ws.Range("E10").Value = '  -6.55%  '
$ws.Range("D11").Value = '0.0810'
$ws.Range("E11").Value = '  -3.55%  '
$ws.Range("D12").Value = '7.70'
$ws.Range("E12").Value = '  -4.56%  '
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").Value = '2.946.89'
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("D15").Value = '2.559.88'
$ws.Range("E15").Value = '  -1.78%  '
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").Value = '14.18'
$ws.Range("D18").Value = '43.150.98'
$ws.Range("E18").Value = '  -6.44%  '
$ws.Range("D19").Value = '13.01'
$ws.Range("E19").Value = '  +3.30%  '
$ws.Range("D20").Value = '0.0₃0981'
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").Value = '72.42'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '260.42'
$ws.Range("E23").Value = '  -10.20%  '
$ws.Range("D24").Value = '2.92'
$ws.Range("E24").Value = '  -3.27%  '
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("E26").Value = '  -4.56%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  -6.33%  '
$ws.Range("D29").Value = '37.25'
$ws.Range("E29").Value = '  -3.66%  '
$ws.Range("D31").Value = '6.01'
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("D32").Value = '154.42'
$ws.Range("E32").Value = '  -3.05%  '
$ws.Range("D33").Value = '2.17'
$ws.Range("E33").Value = '  -2.53%  '
$ws.Range("E35").Value = '  -6.73%  '
$ws.Range("D36").Value = '0.0799'
$ws.Range("E36").Value = '  -4.63%  '
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("D39").Value = '16.71'
$ws.Range("E39").Value = '  +6.59%  '
$ws.Range("D40").Value = '23.27'
$ws.Range("E40").Value = '  +8.85%  '
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("D42").Value = '0.0313'
$ws.Range("E42").Value = '  -4.62%  '
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("D44").Value = '2.071.45'
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '85.90'
$ws.Range("E46").Value = '  -10.40%  '
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("D48").Value = '2.802.93'
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").Value = '8.74'
$ws.Range("E50").Value = '  -6.33%  '
$ws.Range("D51").Value = '104.49'
$ws.Range("E51").Value = '  -4.27%  '
